# Fruta / hortaliza, semanal
# Insert 3 new weekly report rows at the top of the existing data block
# (before current row 289), shifting all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 289 (existing rows 289:346 shift to 292:349)
$ws.Rows("289:291").Insert()

# New row 289: Carson / Especial
$ws.Range("A289").Value = 10
$ws.Range("B289").Value = "Vega Modelo de Temuco"
$ws.Range("C289").Value = "La Araucanía"
$ws.Range("D289").Value = 44943
$ws.Range("E289").Value = 9
$ws.Range("F289").Value = "Fruta"
$ws.Range("G289").Value = 100103
$ws.Range("H289").Value = "Frutos de hueso (carozo)"
$ws.Range("I289").Value = 100103004
$ws.Range("J289").Value = "Durazno"
$ws.Range("K289").Value = "Carson"
$ws.Range("L289").Value = "Especial"
$ws.Range("M289").Value = 5
$ws.Range("N289").Value = 440000
$ws.Range("O289").Value = 440000
$ws.Range("P289").Value = 440000
$ws.Range("Q289").Value = "$/bins (400 kilos)"
$ws.Range("R289").Value = "Región de O'Higgins"
$ws.Range("S289").Value = 1100
$ws.Range("T289").Value = 400

# New row 290: Carson / Primera
$ws.Range("A290").Value = 10
$ws.Range("B290").Value = "Vega Modelo de Temuco"
$ws.Range("C290").Value = "La Araucanía"
$ws.Range("D290").Value = 44943
$ws.Range("E290").Value = 9
$ws.Range("F290").Value = "Fruta"
$ws.Range("G290").Value = 100103
$ws.Range("H290").Value = "Frutos de hueso (carozo)"
$ws.Range("I290").Value = 100103004
$ws.Range("J290").Value = "Durazno"
$ws.Range("K290").Value = "Carson"
$ws.Range("L290").Value = "Primera"
$ws.Range("M290").Value = 8
$ws.Range("N290").Value = 390000
$ws.Range("O290").Value = 390000
$ws.Range("P290").Value = 390000
$ws.Range("Q290").Value = "$/bins (400 kilos)"
$ws.Range("R290").Value = "Región de O'Higgins"
$ws.Range("S290").Value = 975
$ws.Range("T290").Value = 400

# New row 291: Elegant Lady / Primera
$ws.Range("A291").Value = 10
$ws.Range("B291").Value = "Vega Modelo de Temuco"
$ws.Range("C291").Value = "La Araucanía"
$ws.Range("D291").Value = 44943
$ws.Range("E291").Value = 9
$ws.Range("F291").Value = "Fruta"
$ws.Range("G291").Value = 100103
$ws.Range("H291").Value = "Frutos de hueso (carozo)"
$ws.Range("I291").Value = 100103004
$ws.Range("J291").Value = "Durazno"
$ws.Range("K291").Value = "Elegant Lady"
$ws.Range("L291").Value = "Primera"
$ws.Range("M291").Value = 8
$ws.Range("N291").Value = 420000
$ws.Range("O291").Value = 420000
$ws.Range("P291").Value = 420000
$ws.Range("Q291").Value = "$/bins (400 kilos)"
$ws.Range("R291").Value = "Región de O'Higgins"
$ws.Range("S291").Value = 1050
$ws.Range("T291").Value = 400
